$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '24.646.75'
$ws.Range("E2").Value = '  +1.66%  '

# Row 3
$ws.Range("D3").Value = '1.703.44'
$ws.Range("E3").Value = '  +1.18%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.05%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '308.69'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.25%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9979'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.05%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3727'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.59%  '

# Row 8
$ws.Range("B8").Value = 'Cardano'
$ws.Range("C8").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3436'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.56%  '

# Row 9
$ws.Range("B9").Value = 'OKB'
$ws.Range("C9").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '48.75'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.88%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.181'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.34%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07434'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.40%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9982'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.11%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.84'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.23%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.218'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.30%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.917'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.78%  '

# Row 16
$ws.Range("D16").Value = '1.706.31'
$ws.Range("E16").Value = '  +1.58%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001117'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.21%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9983'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.11%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06687'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.59%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '83.20'
$ws.Range("D20").Style = "Normal"

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.07'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.96%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.331'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.05%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.16'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +9.12%  '

# Row 24
$ws.Range("D24").Value = '24.618.20'
$ws.Range("E24").Value = '  +1.62%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.411'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.01%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.755'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.25%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.06'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.03%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '149.53'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.74%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '130.78'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.79%  '

# Row 30
$ws.Range("D30").Value = '1.895.00'
$ws.Range("E30").Value = '  +1.65%  '

# Row 31
$ws.Range("E31").Value = '  +17.00%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.715'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.18%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.192'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.15%  '

# Row 34
$ws.Range("B34").Value = 'Aptos'
$ws.Range("C34").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '13.64'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +8.05%  '

# Row 35
$ws.Range("B35").Value = 'Stellar'
$ws.Range("C35").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.08768'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.10%  '

# Row 36
$ws.Range("B36").Value = 'WEMIXTOKEN'
$ws.Range("C36").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.767'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.50%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.495'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.41%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06511'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.03%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.900'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.83%  '

# Row 40
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.02363'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.15%  '

# Row 41
$ws.Range("B41").Value = 'Algorand'
$ws.Range("C41").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2212'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.05%  '

# Row 42
$ws.Range("E42").Value = '  -1.22%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6379'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.52%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9973'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.05%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.75'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.13%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6057'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.17%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.803'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.33%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.108'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.37%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '128.87'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.89%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07248'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.89%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '78.83'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.72%  '
